$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Remark" column (E) ---
# Set values in the order needed so new shared-string entries are created
# in the same sequence as the target workbook.

# Header
$ws.Range("E2").Value = "Remark"
$ws.Range("E2").Font.Bold = $true

# Remarks for existing rows (first occurrences create new shared strings,
# later ones reuse them)
$ws.Range("E5").Value = "ring"
$ws.Range("E6").Value = "ring"

$ws.Range("E7").Value = "invalid number"

$ws.Range("E16").Value = "ring"

$ws.Range("E23").Value = "invaild"

$ws.Range("E27").Value = "call today evening 8:00 PM"

# --- Correct the region / cluster name on row 27 ---
$ws.Range("C27").Value = "State : Madhya PradeshDistrict : GUNABlock : RAGHOGARHVillage : GOVINDPURACluster : HSS Kumbhraj"

$ws.Range("E30").Value = "not intrested"

$ws.Range("E32").Value = "already software own , but intrested"

$ws.Range("E33").Value = "invalid number"

$ws.Range("E35").Value = "not intrested"

$ws.Range("E36").Value = "incoming call not available"

# --- Split the combined phone numbers in D16, keep only the second number ---
$ws.Range("D16").Value = 9755000536

# --- Restore the active view/selection to match the latest edit ---
$win = $wb.Windows.Item(1)
$win.ScrollRow = 14
$win.ScrollColumn = 3
$ws.Range("E36").Select()
